$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.55986762046814
$ws.Range("B1").Value = 1.970776915550232
$ws.Range("C1").Value = 3.525780916213989
$ws.Range("D1").Value = 1.392233490943909
$ws.Range("E1").Value = 0.8290205001831055
